# chore: add monthly employment outputs
# Updates monitoring-indicator counts across the summary/report sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: 요약_권역별 ---
$ws1 = $wb.Worksheets.Item("요약_권역별")
$ws1.Range("C3").Value = 28
$ws1.Range("E3").Value = 2
$ws1.Range("C4").Value = 19
$ws1.Range("D4").Value = 6
$ws1.Range("D5").Value = 1
$ws1.Range("E5").Value = 0
$ws1.Range("C6").Value = 30
$ws1.Range("D6").Value = 1
$ws1.Range("C7").Value = 28
$ws1.Range("D7").Value = 1
$ws1.Range("C9").Value = 16
$ws1.Range("E9").Value = 1
$ws1.Range("C10").Value = 16
$ws1.Range("D10").Value = 1
$ws1.Range("C13").Value = 13
$ws1.Range("D13").Value = 4

# --- Sheet: 요약_전월대비 ---
$ws2 = $wb.Worksheets.Item("요약_전월대비")
$ws2.Range("C3").Value = 28
$ws2.Range("E3").Value = 2
$ws2.Range("C4").Value = 19
$ws2.Range("D4").Value = 6
$ws2.Range("D5").Value = 1
$ws2.Range("E5").Value = 0
$ws2.Range("C6").Value = 30
$ws2.Range("D6").Value = 1
$ws2.Range("C7").Value = 28
$ws2.Range("D7").Value = 1
$ws2.Range("C9").Value = 16
$ws2.Range("E9").Value = 1
$ws2.Range("C10").Value = 16
$ws2.Range("D10").Value = 1
$ws2.Range("C13").Value = 13
$ws2.Range("D13").Value = 4

# --- Sheet: 3개월연속_시군 ---
$ws4 = $wb.Worksheets.Item("3개월연속_시군")
$ws4.Range("C5").Value = "주의"
$ws4.Range("E5").Value = "관심"
$ws4.Range("C11").Value = "관심"
$ws4.Range("E11").Value = "주의"
$ws4.Range("C14").Value = "관심"
$ws4.Range("E14").Value = "주의"

# --- Sheet: 주요지역_시군 ---
$ws6 = $wb.Worksheets.Item("주요지역_시군")
$ws6.Range("C52").Value = 0
$ws6.Range("D52").Value = 2
$ws6.Range("E52").Value = 2
